$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)

# Resize the textbox to reflect the removal of two list items (EMU -> points)
$shp.Height = 4524315 / 914400 * 72

$tr = $shp.TextFrame.TextRange

# Remove the "Kafka Messaging System, Zoo Keeper" paragraph and the blank
# bullet paragraph that precedes it. Delete from the highest index first so
# the lower index stays valid after the first deletion.
$kafkaPara = $tr.Paragraphs(14, 1)
$kafkaPara.Delete()

$blankPara = $tr.Paragraphs(13, 1)
$blankPara.Delete()
